$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: shift numeric data (columns C:T) for rows 4..29 down by 2 rows.
# Process from the bottom (row 29 -> 31) upward to row 4 -> 6 so source rows are read before being overwritten.
for ($n = 29; $n -ge 4; $n--) {
    $dest = $n + 2
    for ($c = 3; $c -le 20; $c++) {
        $ws.Cells.Item($dest, $c).Value = $ws.Cells.Item($n, $c).Value()
    }
}

# --- Step 2: set the two brand-new rows (4 and 5) with freshly simulated data.
$ws.Cells.Item(4,3).Value = 0.8094219110980454
$ws.Cells.Item(4,4).Value = 0.9382087744045872
$ws.Cells.Item(4,5).Value = 0.8094219110980454
$ws.Cells.Item(4,6).Value = 1.096124913496067
$ws.Cells.Item(4,7).Value = 1.045821774604294
$ws.Cells.Item(4,8).Value = 0.8466941654706673
$ws.Cells.Item(4,9).Value = 1.44287730690868
$ws.Cells.Item(4,10).Value = 0.9382087744045872
$ws.Cells.Item(4,11).Value = 0.9382087744045872
$ws.Cells.Item(4,12).Value = 1.096124913496067
$ws.Cells.Item(4,13).Value = 0.9527734122970564
$ws.Cells.Item(4,14).Value = 0.9527734122970564
$ws.Cells.Item(4,15).Value = 0.9174136633549267
$ws.Cells.Item(4,16).Value = 0.9479185329995666
$ws.Cells.Item(4,17).Value = 0.9479185329995666
$ws.Cells.Item(4,18).Value = 0.9454910933508217
$ws.Cells.Item(4,19).Value = 0.9454910933508217
$ws.Cells.Item(4,20).Value = 1.029858140997057
$ws.Cells.Item(5,3).Value = 1.283930090580594
$ws.Cells.Item(5,4).Value = 0.789831345357614
$ws.Cells.Item(5,5).Value = 1.283930090580594
$ws.Cells.Item(5,6).Value = 1.119805599905353
$ws.Cells.Item(5,7).Value = 0.9960066618511904
$ws.Cells.Item(5,8).Value = 0.8487985641260259
$ws.Cells.Item(5,9).Value = 1.669255030527202
$ws.Cells.Item(5,10).Value = 0.789831345357614
$ws.Cells.Item(5,11).Value = 0.789831345357614
$ws.Cells.Item(5,12).Value = 1.119805599905353
$ws.Cells.Item(5,13).Value = 1.201867845242973
$ws.Cells.Item(5,14).Value = 1.201867845242973
$ws.Cells.Item(5,15).Value = 1.084178084870657
$ws.Cells.Item(5,16).Value = 1.064522345281187
$ws.Cells.Item(5,17).Value = 1.064522345281187
$ws.Cells.Item(5,18).Value = 0.9958495953002937
$ws.Cells.Item(5,19).Value = 0.9958495953002937
$ws.Cells.Item(5,20).Value = 1.117937882057997

# --- Step 3: set/refresh the HKL label text in column B for rows 4..31 (labels shift along with the data; "Thomas Hex" becomes "Matthies Hex").
$ws.Cells.Item(4,2).Value = "Holden"
$ws.Cells.Item(5,2).Value = "Rizzie Spiral"
$ws.Cells.Item(6,2).Value = "RotRing OmegaMax-90"
$ws.Cells.Item(7,2).Value = "Equal Angle"
$ws.Cells.Item(8,2).Value = "Tilt Rotate"
$ws.Cells.Item(9,2).Value = "CLR"
$ws.Cells.Item(10,2).Value = "Rizzie Hex"
$ws.Cells.Item(11,2).Value = "Matthies Hex"
$ws.Cells.Item(12,2).Value = "Tilt Rotate_Partial"
$ws.Cells.Item(13,2).Value = "RotRing OmegaMax-60"
$ws.Cells.Item(14,2).Value = "Equal Angle_Partial"
$ws.Cells.Item(15,2).Value = "Rizzie Hex_Partial"
$ws.Cells.Item(16,2).Value = "ND Single"
$ws.Cells.Item(17,2).Value = "RD Single"
$ws.Cells.Item(18,2).Value = "TD Single"
$ws.Cells.Item(19,2).Value = "Morris Single"
$ws.Cells.Item(20,2).Value = "Ring Perpendicular to ND"
$ws.Cells.Item(21,2).Value = "Ring Perpendicular to RD"
$ws.Cells.Item(22,2).Value = "Ring Perpendicular to TD"
$ws.Cells.Item(23,2).Value = "OffsetFTD"
$ws.Cells.Item(24,2).Value = "OffsetATD"
$ws.Cells.Item(25,2).Value = "OffsetF45"
$ws.Cells.Item(26,2).Value = "OffsetA45"
$ws.Cells.Item(27,2).Value = "OffsetFRD"
$ws.Cells.Item(28,2).Value = "OffsetARD"
$ws.Cells.Item(29,2).Value = "Gaussian Quadrature"
$ws.Cells.Item(30,2).Value = "Michael-CCHex"
$ws.Cells.Item(31,2).Value = "Michael-SNHex"

# --- Step 4: set column A (HKL index) and give it the same style as existing indexed cells for the two brand-new rows (30, 31).
$ws.Cells.Item(30,1).Value = 28
$ws.Cells.Item(31,1).Value = 29
$ws.Cells.Item(30,1).Style = $ws.Cells.Item(29,1).Style
$ws.Cells.Item(31,1).Style = $ws.Cells.Item(29,1).Style
